# Update countries & provincias Spain
# Applies the 27-May-2020 23:05 data refresh to the "Pais" sheet:
#  - refreshed case counts for a handful of countries (including the
#    sort-order-driving columns), which in turn re-shuffles a few
#    neighbouring rows' country names because the sheet is kept sorted
#    by total cases descending
#  - bumps the "Datos actualizados" timestamp in the title cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Title / timestamp cell
$ws.Range("A1").Value = "Datos actualizados a 27 de Mayo de 2020 a las 23:05"

function Set-CountryRow {
    param([int]$RowNum, [string]$Country, [double]$B, [double]$C, [double]$D, [double]$E, [double]$F, [double]$G, [double]$H)
    $ws.Range("A$RowNum").Value = $Country
    $ws.Range("B$RowNum").Value = $B
    $ws.Range("C$RowNum").Value = $C
    $ws.Range("D$RowNum").Value = $D
    $ws.Range("E$RowNum").Value = $E
    $ws.Range("F$RowNum").Value = $F
    $ws.Range("G$RowNum").Value = $G
    $ws.Range("H$RowNum").Value = $H
}

# Row 4 - Estados Unidos (refreshed counts, country unchanged)
Set-CountryRow 4 "Estados Unidos" 1740694 15419 483866 1155015 0 1241 101813

# Row 11 - Alemania (refreshed counts, country unchanged)
Set-CountryRow 11 "Alemania" 181895 607 162800 10562 0 35 8533

# Rows 151-154: Mauritania (unchanged, row 150) is followed by Uganda now
# instead of Suazilandia/Liberia/Yemen/Uganda
Set-CountryRow 151 "Uganda"       281 28 69  212 0 0 0
Set-CountryRow 152 "Suazilandia"  272 11 168 102 0 0 2
Set-CountryRow 153 "Liberia"      266 0  144 95  0 1 27
Set-CountryRow 154 "Yemen"        256 7  10  193 0 4 53

# Rows 174-178: Bahamas (unchanged, row 173) is followed by Libia now
# instead of Monaco/Barbados/Comoras/Liechtenstein
Set-CountryRow 174 "Libia"         99 22 40 55 0 1 4
Set-CountryRow 175 "Monaco"        98 0  90 4  0 0 4
Set-CountryRow 176 "Barbados"      92 0  76 9  0 0 7
Set-CountryRow 177 "Comoras"       87 0  24 61 0 1 2
Set-CountryRow 178 "Liechtenstein" 82 0  55 26 0 0 1

# Rows 199-201: Fiyi (unchanged, row 198) is followed by Santa Lucia /
# Nueva Caledonia / Belice instead of Belice / Nueva Caledonia / Santa Lucia
Set-CountryRow 199 "Nueva Caledonia" 18 0 18 0 0 0 0
Set-CountryRow 200 "Santa Lucia"     18 0 18 0 0 0 0
Set-CountryRow 201 "Belice"          18 0 16 0 0 0 2
